# "Looking at CLalpha in radians"
# Adds an "alpha(rad)" helper column (M) that converts the angle-of-attack
# column (A) from degrees to radians, and a "CLALPHA" column (O) that
# estimates the lift-curve slope from the Cl data using that radian column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), same header style as the existing headers ---
$ws.Range("M1").Value = "alpha(rad)"
$ws.Range("O1").Value = "CLALPHA"
$ws.Range("M1").VerticalAlignment = -4108
$ws.Range("O1").VerticalAlignment = -4108

# --- alpha(rad) = alpha(deg) * (PI()/180), filled down rows 2:42 ---
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 13).Formula = "=A" + $r + "*(PI()/180)"
}

# --- CLALPHA: lift-curve slope estimated from row 2 and row 23 (linear region) ---
$ws.Range("O2").Formula = "=(B23-B2)/(M23-M2)"

# --- Column M width to roughly match the autofit width used for the new header ---
$ws.Columns.Item(13).ColumnWidth = 9.3

# --- Keep the active selection where the author left it: O3 ---
[void]$ws.Range("O3").Select()

# --- The source file also dropped a stale external-workbook reference
#     (left over from a formula that no longer exists) when it was re-saved.
#     Try to break it if the host bridges the call; harmless if not. ---
try {
    $names = $wb.LinkSources()
    if ($names) {
        foreach ($n in $names) {
            $wb.BreakLink($n, 1)
        }
    }
} catch {
}
